$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-9
# from serial 45221 (2023-10-22) to serial 45224 (2023-10-25)
for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45221) {
        $cell.Value2 = 45224
    }
}
